$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The pre-existing Month/Year cells (D2, D6, D21, D22 - *before* the new row
# is inserted below) currently hold real Excel date serial numbers. They
# need to hold the literal text of the formatted date instead, while keeping
# their existing date number-format style. Writing the text as a formula
# that evaluates to the literal string, then collapsing the formula down to
# its resulting value via copy / paste-special-values, yields a plain shared
# string cell (t="s") without reapplying / creating any new number format,
# so styles.xml stays untouched.
# (Doing this before the row insert below also keeps the shared-string
# table ordering identical to the authored edit.)
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 4).Formula = "=""2025 July"""
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)

$ws.Cells.Item(6, 4).Formula = "=""2024 December 14"""
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)

$ws.Cells.Item(21, 4).Formula = "=""2018 April"""
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)

$ws.Cells.Item(22, 4).Formula = "=""2018 March"""
$ws.Cells.Item(22, 4).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Insert a brand new row 2 (Akriti Kumari / IIT Hyderabad award) and push all
# existing rows (old row 2 .. old row 23) down by one.
# ---------------------------------------------------------------------------
$ws.Rows("2:2").Insert()

$ws.Cells.Item(2, 1).Value = "Akriti Kumari "
$ws.Cells.Item(2, 2).Value = "BUILD (Bold and Unique Idea Led Development) project  to develop a diagnostic kit for Urinary Tract Infections (UTI) and Antimicrobial Resistance (AMR)."
$ws.Cells.Item(2, 3).Value = "IIT Hyderabad"
$ws.Cells.Item(2, 5).Value = "Dr. Gunjan Mehta"

# Column D ("Month/Year") here also needs to hold plain text "2025 October"
# rather than being auto-recognised as a date by Excel - same trick as above.
$ws.Cells.Item(2, 4).Formula = "=""2025 October"""
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column E now needs an explicit width, matching the widened "Awarded By"
# column seen after the edit (45.85546875 characters in the saved OOXML;
# the nearest value reachable through the ColumnWidth property is used).
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 45

# ---------------------------------------------------------------------------
# Restore the final selected cell shown in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("E2").Select() | Out-Null
